$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths for A and B (target stored width 11.3046875; engine quantizes
# ColumnWidth to a 1/6-character grid internally, so 10.5 is the closest
# input that lands nearest the target stored width)
$ws.Columns.Item(1).ColumnWidth = 10.5
$ws.Columns.Item(2).ColumnWidth = 10.5

# Row 2
$ws.Range("A2").Value = 45644.148448182874
$ws.Range("B2").Value = 45652.0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = $true
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 5

# Row 3
$ws.Range("A3").Value = 45644.14986158565
$ws.Range("B3").Value = 45650.149886782405
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = $true
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 9

# Number formats for date columns
$ws.Range("A2:B3").NumberFormat = "yyyy-mm-dd"
